# DS 775 Wk 10 - add problem 20.6-2 worksheet and refresh the frozen
# random draws on 20.1-1 / 20.1-2 that feed the earlier problems.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 20.1-1 : replace the live random draws in C3:C8 with the values that
# were sampled when the workbook was last "frozen" (paste-values).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("20.1-1")

$ws1.Range("C3").Value = 0.3039
$ws1.Range("C4").Value = 0.7914
$ws1.Range("C5").Value = 0.8543
$ws1.Range("C6").Value = 0.6902
$ws1.Range("C7").Value = 0.3004
$ws1.Range("C8").Value = 0.0383

$ws1.Activate()
$ws1.Range("G17").Select()

# ---------------------------------------------------------------------
# 20.1-2 : same idea for the second weather table, C17:C26.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("20.1-2")

$ws2.Range("C17").Value = 0.3004
$ws2.Range("C18").Value = 0.0383
$ws2.Range("C19").Value = 0.3883
$ws2.Range("C20").Value = 0.6052
$ws2.Range("C21").Value = 0.2231
$ws2.Range("C22").Value = 0.425
$ws2.Range("C23").Value = 0.3729
$ws2.Range("C24").Value = 0.7983
$ws2.Range("C25").Value = 0.234
$ws2.Range("C26").Value = 0.0082

$ws2.Activate()
$ws2.Range("G18").Select()

# ---------------------------------------------------------------------
# Add the new worksheet for problem 20.6-2 at the end of the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "20.6-2"

$ws5.Range("A1").Value = "20.6-2"
$ws5.Range("A1").Style = "Heading 1"

$ws5.Range("C2").Value = "mean"
$ws5.Range("D2").Value = "sd"
$ws5.Range("E2").Value = "distrbution"
$ws5.Range("F2").Value = "yr"

$ws5.Range("B3").Value = "cost to purchase"
$ws5.Range("C3").Value = 1000
$ws5.Range("F3").Value = 0

$ws5.Range("B4").Value = "construction"
$ws5.Range("C4").Value = 2000
$ws5.Range("D4").Value = 0.2
$ws5.Range("E4").Value = "triangle"
$ws5.Range("F4").Value = 1

$ws5.Range("B5").Value = "annual profit"
$ws5.Range("C5").Value = 700
$ws5.Range("D5").Value = 700
$ws5.Range("E5").Value = "normal"
$ws5.Range("F5").Value = "2,3,4,5"

$ws5.Range("B6").Value = "sell"
$ws5.Range("C6").Value = "4000-8000"
$ws5.Range("F6").Value = 5

$ws5.Range("I9").Value = "End of Year "
$ws5.Range("J9").Value = "In"
$ws5.Range("K9").Value = "Out"
$ws5.Range("L9").Value = "Net"

$ws5.Range("I10").Value = 0
$ws5.Range("K10").Value = -1000
$ws5.Range("L10").Formula = "=J10+K10"

$ws5.Range("I11").Value = 1
$ws5.Range("K11").Value = -2000
$ws5.Range("K11").Interior.Color = 65535
$ws5.Range("L11").Formula = "=J11+K11"

$ws5.Range("I12").Value = 2
$ws5.Range("J12").Value = 700
$ws5.Range("J12").Interior.Color = 65535
$ws5.Range("L12").Formula = "=J12+K12"

$ws5.Range("I13").Value = 3
$ws5.Range("J13").Value = 700
$ws5.Range("J13").Interior.Color = 65535
$ws5.Range("L13").Formula = "=J13+K13"

$ws5.Range("I14").Value = 4
$ws5.Range("J14").Value = 700
$ws5.Range("J14").Interior.Color = 65535
$ws5.Range("L14").Formula = "=J14+K14"

$ws5.Range("I15").Value = 5
$ws5.Range("J15").Value = 700
$ws5.Range("J15").Interior.Color = 65535
$ws5.Range("L15").Formula = "=J15+K15"

$ws5.Range("L16").Formula = "=SUM(L10:L15)"

$ws5.Columns.Item(1).AutoFit() | Out-Null
$ws5.Columns.Item(2).AutoFit() | Out-Null
$ws5.Columns.Item(9).AutoFit() | Out-Null

$ws5.Activate()
$ws5.Range("L16").Select()
